# Apply targeted cell value updates to Sheet1, rows 5 and 11
# (odds data refresh for Jogos_da_Semana_FlashScore_2024-10-24.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    # Row 5
    "G5"  = 1.6
    "I5"  = 6.25
    "J5"  = 2.3
    "K5"  = 2
    "L5"  = 7
    "M5"  = 1.08
    "N5"  = 8
    "O5"  = 1.5
    "P5"  = 2.5
    "Q5"  = 2.5
    "R5"  = 1.5
    "U5"  = 2.63
    "V5"  = 1.44
    "W5"  = 4.75
    "X5"  = 6
    "Y5"  = 9.5
    "Z5"  = 11
    "AA5" = 17
    "AC5" = 6.5
    "AD5" = 7.5
    "AE5" = 26
    "AF5" = 126
    "AH5" = 29
    "AI5" = 21
    "AK5" = 51
    "AL5" = 67
    "AN5" = 3.25
    "AO5" = 9
    "AS5" = 301
    "AW5" = 7.5
    "BA5" = 251

    # Row 11
    "H11"  = 3.55
    "J11"  = 3.9
    "K11"  = 2.37
    "L11"  = 2.25
    "P11"  = 3.4
    "S11"  = 1.33
    "T11"  = 3.14
    "W11"  = 12.5
    "AC11" = 11.75
    "AD11" = 7
    "AG11" = 8
    "AH11" = 9.5
    "AJ11" = 16
    "AL11" = 23
    "AN11" = 6.1
    "AO11" = 19
    "AP11" = 20
    "AQ11" = 80
    "AR11" = 90
    "AS11" = 175
    "AT11" = 3.55
    "AU11" = 6.3
    "AV11" = 40
    "AW11" = 4.05
    "AX11" = 8.5
    "AY11" = 13.5
    "AZ11" = 26
    "BA11" = 40
    "BB11" = 120
    "BC11" = 500
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
